$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Broaden the "Nur in Fragenpool" (question-pool-only) field hints to
# "Optional, nur in Fragenpool" across the whole sheet, and rename the
# generic "Fill in the blanks" question-type label to include the German
# qualifier.
$ws.Cells.Replace("Nur in Fragenpool", "Optional, nur in Fragenpool")
$ws.Cells.Replace("Fragetpy: Fill in the blanks", "Fragetpy: Fill in the blanks (Lückentext)")

# Context help: put the selection on the first optional/context-help cell.
$ws.Range("D4").Select()
